$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Update "Ultima actualizacion" timestamp on all three sheets
$ws1.Range("A2").Value = "Última actualización: 17:56:03"
$ws2.Range("A2").Value = "Última actualización: 17:56:03"
$ws3.Range("A2").Value = "Última actualización: 17:56:03"

# Update total row count on sheet1
$ws1.Range("A3").Value = "Total filas: 299"

# Apply cell-level corrections to existing data rows (re-sorted / re-scraped entries)
$ws1.Cells.Item(23,1).Value = "06:17:28"
$ws1.Cells.Item(23,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(23,4).Value = 64

$ws1.Cells.Item(24,1).Value = "05:57:13"
$ws1.Cells.Item(24,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(24,4).Value = 84

$ws1.Cells.Item(33,1).Value = "06:46:50"
$ws1.Cells.Item(33,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(33,4).Value = 74

$ws1.Cells.Item(34,1).Value = "06:17:28"
$ws1.Cells.Item(34,3).Value = "17_ROMERO"
$ws1.Cells.Item(34,4).Value = 103

$ws1.Cells.Item(106,1).Value = "11:47:17"
$ws1.Cells.Item(106,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(106,4).Value = 5

$ws1.Cells.Item(108,1).Value = "11:52:01"
$ws1.Cells.Item(108,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(108,4).Value = 0

$ws1.Cells.Item(135,1).Value = "11:34:59"
$ws1.Cells.Item(135,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(135,4).Value = 62

$ws1.Cells.Item(136,1).Value = "10:50:41"
$ws1.Cells.Item(136,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(136,4).Value = 106

$ws1.Cells.Item(142,1).Value = "10:50:41"
$ws1.Cells.Item(142,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(142,4).Value = 118

$ws1.Cells.Item(143,1).Value = "11:47:17"
$ws1.Cells.Item(143,3).Value = "14_ABASTO"
$ws1.Cells.Item(143,4).Value = 61

$ws1.Cells.Item(144,1).Value = "11:11:33"
$ws1.Cells.Item(144,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(144,4).Value = 97

$ws1.Cells.Item(209,1).Value = "14:53:07"
$ws1.Cells.Item(209,3).Value = "10_OLMOS"
$ws1.Cells.Item(209,4).Value = 60

$ws1.Cells.Item(210,1).Value = "13:56:11"
$ws1.Cells.Item(210,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(210,4).Value = 117

$ws1.Cells.Item(211,3).Value = "15X38_ABASTO"

$ws1.Cells.Item(259,1).Value = "17:48:33"
$ws1.Cells.Item(259,3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(259,4).Value = 0

$ws1.Cells.Item(260,1).Value = "16:44:12"
$ws1.Cells.Item(260,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(260,4).Value = 64

$ws1.Cells.Item(261,3).Value = "215B_EL PATO"

$ws1.Cells.Item(267,1).Value = "17:56:03"
$ws1.Cells.Item(267,2).Value = "17:56"
$ws1.Cells.Item(267,3).Value = "10_OLMOS"
$ws1.Cells.Item(267,4).Value = 0

$ws1.Cells.Item(268,1).Value = "17:56:03"
$ws1.Cells.Item(268,2).Value = "17:57"
$ws1.Cells.Item(268,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(268,4).Value = 1

$ws1.Cells.Item(269,1).Value = "16:28:03"
$ws1.Cells.Item(269,2).Value = "17:58"
$ws1.Cells.Item(269,3).Value = "17_ROMERO"
$ws1.Cells.Item(269,4).Value = 90

$ws1.Cells.Item(270,1).Value = "16:14:44"
$ws1.Cells.Item(270,2).Value = "18:00"
$ws1.Cells.Item(270,3).Value = "10_OLMOS"
$ws1.Cells.Item(270,4).Value = 106

$ws1.Cells.Item(271,2).Value = "18:05"
$ws1.Cells.Item(271,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(271,4).Value = 111

$ws1.Cells.Item(272,1).Value = "16:28:03"
$ws1.Cells.Item(272,2).Value = "18:06"
$ws1.Cells.Item(272,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(272,4).Value = 98

$ws1.Cells.Item(273,1).Value = "16:14:44"
$ws1.Cells.Item(273,2).Value = "18:10"
$ws1.Cells.Item(273,3).Value = "15_ABASTO"
$ws1.Cells.Item(273,4).Value = 116

$ws1.Cells.Item(274,1).Value = "16:14:44"
$ws1.Cells.Item(274,2).Value = "18:10"
$ws1.Cells.Item(274,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(274,4).Value = 116

$ws1.Cells.Item(275,2).Value = "18:17"
$ws1.Cells.Item(275,3).Value = "10_OLMOS"
$ws1.Cells.Item(275,4).Value = 109

$ws1.Cells.Item(276,1).Value = "16:52:42"
$ws1.Cells.Item(276,2).Value = "18:21"
$ws1.Cells.Item(276,3).Value = "215C_EL PATO"
$ws1.Cells.Item(276,4).Value = 89

$ws1.Cells.Item(277,1).Value = "16:28:03"
$ws1.Cells.Item(277,2).Value = "18:22"
$ws1.Cells.Item(277,3).Value = "215C_EL PATO"
$ws1.Cells.Item(277,4).Value = 114

$ws1.Cells.Item(278,1).Value = "16:28:03"
$ws1.Cells.Item(278,2).Value = "18:25"
$ws1.Cells.Item(278,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(278,4).Value = 117

$ws1.Cells.Item(279,1).Value = "16:52:42"
$ws1.Cells.Item(279,2).Value = "18:29"
$ws1.Cells.Item(279,3).Value = "14_ABASTO"
$ws1.Cells.Item(279,4).Value = 97

$ws1.Cells.Item(280,2).Value = "18:30"
$ws1.Cells.Item(280,3).Value = "14_ABASTO"
$ws1.Cells.Item(280,4).Value = 113

$ws1.Cells.Item(281,1).Value = "17:48:33"
$ws1.Cells.Item(281,2).Value = "18:32"
$ws1.Cells.Item(281,4).Value = 44

$ws1.Cells.Item(282,1).Value = "16:37:06"
$ws1.Cells.Item(282,2).Value = "18:36"
$ws1.Cells.Item(282,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(282,4).Value = 119

$ws1.Cells.Item(283,2).Value = "18:36"
$ws1.Cells.Item(283,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(283,4).Value = 83

$ws1.Cells.Item(284,1).Value = "17:36:10"
$ws1.Cells.Item(284,2).Value = "18:37"
$ws1.Cells.Item(284,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(284,4).Value = 61

$ws1.Cells.Item(285,2).Value = "18:41"
$ws1.Cells.Item(285,3).Value = "10_OLMOS"
$ws1.Cells.Item(285,4).Value = 88

$ws1.Cells.Item(286,1).Value = "16:52:42"
$ws1.Cells.Item(286,2).Value = "18:45"
$ws1.Cells.Item(286,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(286,4).Value = 113

$ws1.Cells.Item(287,2).Value = "18:52"
$ws1.Cells.Item(287,3).Value = "17_ROMERO"
$ws1.Cells.Item(287,4).Value = 99

$ws1.Cells.Item(288,1).Value = "17:13:39"
$ws1.Cells.Item(288,2).Value = "18:57"
$ws1.Cells.Item(288,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(288,4).Value = 104

$ws1.Cells.Item(289,2).Value = "18:59"
$ws1.Cells.Item(289,3).Value = "14_ABASTO"
$ws1.Cells.Item(289,4).Value = 106

$ws1.Cells.Item(290,2).Value = "19:00"
$ws1.Cells.Item(290,3).Value = "14_ABASTO"
$ws1.Cells.Item(290,4).Value = 84

$ws1.Cells.Item(291,1).Value = "17:13:39"
$ws1.Cells.Item(291,2).Value = "19:03"
$ws1.Cells.Item(291,3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(291,4).Value = 110

$ws1.Cells.Item(292,2).Value = "19:04"
$ws1.Cells.Item(292,3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(292,4).Value = 88

$ws1.Cells.Item(293,1).Value = "17:56:03"
$ws1.Cells.Item(293,2).Value = "19:10"
$ws1.Cells.Item(293,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(293,4).Value = 74

$ws1.Cells.Item(294,1).Value = "17:48:33"
$ws1.Cells.Item(294,2).Value = "19:12"
$ws1.Cells.Item(294,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(294,4).Value = 84

$ws1.Cells.Item(295,1).Value = "17:56:03"
$ws1.Cells.Item(295,2).Value = "19:16"
$ws1.Cells.Item(295,3).Value = "17_ROMERO"
$ws1.Cells.Item(295,4).Value = 80

# Append new scraped rows at the end of the table
$ws1.Cells.Item(296,1).Value = "17:36:10"
$ws1.Cells.Item(296,2).Value = "19:17"
$ws1.Cells.Item(296,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(296,4).Value = 101
$ws1.Cells.Item(296,5).Value = "LP1912"

$ws1.Cells.Item(297,1).Value = "17:36:10"
$ws1.Cells.Item(297,2).Value = "19:17"
$ws1.Cells.Item(297,3).Value = "14X44_ABASTO"
$ws1.Cells.Item(297,4).Value = 101
$ws1.Cells.Item(297,5).Value = "LP1912"

$ws1.Cells.Item(298,1).Value = "17:56:03"
$ws1.Cells.Item(298,2).Value = "19:21"
$ws1.Cells.Item(298,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(298,4).Value = 85
$ws1.Cells.Item(298,5).Value = "LP1912"

$ws1.Cells.Item(299,1).Value = "17:36:10"
$ws1.Cells.Item(299,2).Value = "19:28"
$ws1.Cells.Item(299,3).Value = "215C_EL PATO"
$ws1.Cells.Item(299,4).Value = 112
$ws1.Cells.Item(299,5).Value = "LP1912"

$ws1.Cells.Item(300,1).Value = "17:48:33"
$ws1.Cells.Item(300,2).Value = "19:35"
$ws1.Cells.Item(300,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(300,4).Value = 107
$ws1.Cells.Item(300,5).Value = "LP1912"

$ws1.Cells.Item(301,1).Value = "17:56:03"
$ws1.Cells.Item(301,2).Value = "19:36"
$ws1.Cells.Item(301,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(301,4).Value = 100
$ws1.Cells.Item(301,5).Value = "LP1912"

$ws1.Cells.Item(302,1).Value = "17:48:33"
$ws1.Cells.Item(302,2).Value = "19:39"
$ws1.Cells.Item(302,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(302,4).Value = 111
$ws1.Cells.Item(302,5).Value = "LP1912"

$ws1.Cells.Item(303,1).Value = "17:56:03"
$ws1.Cells.Item(303,2).Value = "19:52"
$ws1.Cells.Item(303,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(303,4).Value = 116
$ws1.Cells.Item(303,5).Value = "LP1912"

$ws1.Cells.Item(304,1).Value = "17:56:03"
$ws1.Cells.Item(304,2).Value = "19:53"
$ws1.Cells.Item(304,3).Value = "225_GOMEZ"
$ws1.Cells.Item(304,4).Value = 117
$ws1.Cells.Item(304,5).Value = "LP1912"
